$d = $word.ActiveDocument

# The first two paragraphs ("Story of a Legacy" title + "By Dorothy Day" byline)
# are being rebuilt as a pandoc-style title block: a Title-styled heading with
# the title words as separate runs, followed by an Authors-styled paragraph
# with the author's name as separate runs (no bold). The old bookmark that
# wrapped the heading paragraph is dropped entirely.

$titlePara = $d.Paragraphs.Item(1)
$authorPara = $d.Paragraphs.Item(2)

$start = $titlePara.Range.Start
$end = $authorPara.Range.End
$target = $d.Range($start, $end)

$blockXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:pStyle w:val="Title"/></w:pPr>
<w:r><w:t xml:space="preserve">Story</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve">of</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve">a</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve">Legacy</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="Authors"/></w:pPr>
<w:r><w:t xml:space="preserve">Dorothy</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve">Day</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$target.InsertXML($blockXml)
